# wc_16.xlsx edit: hide Sheet1/Sheet3, rename Sheet4 -> group_stages,
# re-enter the L/M/N helper-column formulas as range formulas (Excel then
# stores them as shared formulas), clear the one-off VLOOKUP in Sheet2!J6
# to a literal value, and restore the recorded cell selections.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Re-enter formulas for the L/M/N columns so Excel groups them into
#     shared formulas (t="shared") instead of one <f> per cell. ---
$ws2.Range("L2:L9").Formula = "=E2-F2"
$ws2.Range("M2:M9").Formula = "=(E2-F2)/2"
$ws2.Range("N2:N9").Formula = "=C2-D2"

$ws4.Range("L2:L49").Formula = "=E2-F2"
$ws4.Range("M2:M49").Formula = "=(E2-F2)/2"
$ws4.Range("N2:N49").Formula = "=C2-D2"

# --- Sheet2!J6 loses its VLOOKUP formula and becomes a plain literal ---
$ws2.Range("J6").Value = 4

# --- Selections recorded in the saved view state ---
$ws4.Select() | Out-Null
$ws4.Range("W34").Select() | Out-Null

$ws2.Select() | Out-Null
$ws2.Range("J7").Select() | Out-Null

# --- Rename Sheet4 and hide Sheet1 / Sheet3 ---
$ws4.Name = "group_stages"
$wb.Worksheets.Item("Sheet1").Visible = $false
$wb.Worksheets.Item("Sheet3").Visible = $false
